# This workbook holds one weekly "snapshot" row per record (rows 2..41),
# with columns: D=Fecha, M=Volumen, N=Precio minimo, O=Precio maximo,
# P=Precio promedio ponderado, S=Precio $/Kg.
# The commit re-shuffles which date's price data appears on each row
# (a weekly re-sync of the underlying logica_diaria dataset), while every
# other column (A,B,C,E..L,Q,R,T) is identical on every row and stays put.
#
# For target row N (2..41) the new D/M/N/O/P/S values are exactly the
# OLD D/M/N/O/P/S values that used to sit on row SourceMap[N].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 41

# SourceMap[i] = original row number whose (D,M,N,O,P,S) tuple becomes the
# new tuple for target row (firstRow + i), i = 0 .. (lastRow-firstRow)
$SourceMap = @(21,27,12,41,7,19,15,28,6,22,4,2,38,30,40,16,34,31,13,18,37,29,25,23,20,36,35,17,10,26,5,24,8,9,33,11,3,14,39,32)

# Columns (by index) that move together as a tuple for each record.
$cols = @(4,13,14,15,16,19)   # D, M, N, O, P, S

# 1) Snapshot all the old values before writing anything, so that later
#    writes don't clobber values still needed as a source for later rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back the permuted values.
for ($i = 0; $i -lt $SourceMap.Length; $i++) {
    $targetRow = $firstRow + $i
    $sourceRow = $SourceMap[$i]
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value2 = $srcVals[$c]
    }
}
